$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -12.57160000000001
$ws.Range("A3").Value = -22.0717
$ws.Range("E3").Value = 16.0856
$ws.Range("E12").Value = 17.31620000000001
$ws.Range("A14").Value = -21.72070000000001
$ws.Range("A16").Value = -21.78720000000001
$ws.Range("C18").Value = -11.4856
$ws.Range("A21").Value = -20.41609999999998
$ws.Range("A23").Value = -20.64139999999998
$ws.Range("C24").Value = -12.69529999999999
$ws.Range("E24").Value = 16.5429
$ws.Range("A25").Value = -21.92889999999999
$ws.Range("C25").Value = -12.9085
$ws.Range("E25").Value = 16.69750000000001
$ws.Range("A26").Value = -21.10629999999997
$ws.Range("C27").Value = -12.4343
$ws.Range("A29").Value = -21.01189999999998
$ws.Range("C30").Value = -13.3694
$ws.Range("C31").Value = -13.2584
$ws.Range("C39").Value = -12.98470000000001
$ws.Range("A40").Value = -20.14309999999999
$ws.Range("E41").Value = 16.5642
$ws.Range("C42").Value = -13.02839999999999
$ws.Range("C48").Value = -11.4576
$ws.Range("E50").Value = 16.2797
$ws.Range("C51").Value = -11.3613
$ws.Range("C52").Value = -11.16079999999999
$ws.Range("A53").Value = -21.91569999999999
$ws.Range("E53").Value = 16.71860000000002
$ws.Range("C55").Value = -13.9257
$ws.Range("C56").Value = -12.82879999999999
$ws.Range("E56").Value = 16.1237
$ws.Range("A57").Value = -22.01719999999999
$ws.Range("C57").Value = -12.91919999999999
$ws.Range("E57").Value = 16.622
$ws.Range("E58").Value = 16.24620000000001
$ws.Range("A59").Value = -22.0722
$ws.Range("C60").Value = -13.45219999999999
$ws.Range("E61").Value = 16.5325
$ws.Range("E63").Value = 17.35680000000002
$ws.Range("E64").Value = 17.3659
$ws.Range("A65").Value = -21.86889999999999
$ws.Range("A69").Value = -21.60649999999999
$ws.Range("E70").Value = 17.24350000000001
$ws.Range("E72").Value = 16.9466
$ws.Range("C73").Value = -12.68800000000001
$ws.Range("C74").Value = -12.1871
$ws.Range("A79").Value = -20.5804
$ws.Range("A83").Value = -21.889
$ws.Range("E86").Value = 16.5192
$ws.Range("C89").Value = -10.2077
$ws.Range("E89").Value = 17.37950000000003
$ws.Range("C90").Value = -12.8012
$ws.Range("A91").Value = -21.44730000000001
$ws.Range("C92").Value = -10.59829999999999
$ws.Range("A93").Value = -21.17439999999999
$ws.Range("E98").Value = 15.88910000000001
$ws.Range("A100").Value = -21.94829999999999
$ws.Range("E100").Value = 16.3604
$ws.Range("E102").Value = 16.48759999999999
